# Updates the cryptos price-list worksheet with the latest snapshot values
# (commit: "Updated cryptos list on Mon Jul 17 11:19:43 UTC 2023 with GitHub Actions").
#
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). Rows 2-51 hold one coin each.
# Some rows also had their coin swapped with the neighbouring row (the scraper's
# ranking shuffled slightly between runs), so B/C get rewritten there too.
#
# D column strings that look like plain numbers (e.g. "1.003") get coerced to
# a numeric type by Excel's normal text/number inference when just setting
# .Value, which would drop the original text-cell representation used by the
# source data (multi-dot prices like "30.313.52" are fine as-is since Excel
# can't parse them as numbers). For those cells we briefly force a Text
# number format, assign the value, then clear the format again so the cell
# keeps its default (unstyled) appearance but stays a text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.313.52'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").Value = '1.921.23'
$ws.Range("E3").Value = '  -0.44%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7435'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.60'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.44%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.002'
$ws.Range("D7").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '27.43'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.73%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3138'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.75%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06984'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.81%  '
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7750'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.86%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07984'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.48%  '
$ws.Range("D13").Value = '1.926.58'
$ws.Range("E13").Value = '  -0.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.307'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.87'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.92%  '
$ws.Range("D16").Value = '30.290.22'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.29'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '247.00'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.869'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007876'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.39%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.172.50'
$ws.Range("E21").Value = '  -0.54%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.002'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.002'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.650'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.442'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.52'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.02'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.35%  '
$ws.Range("E28").Value = '  -4.31%  '
$ws.Range("E29").Value = '  -5.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.367'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.51%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.546'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.71%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.358'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.080'
$ws.Range("D33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05201'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.69%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.304'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.77%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7497'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.772'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01944'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.787'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.438'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '76.02'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4475'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.950'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.98%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8423'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.40%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.678'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.52'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.836'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.47%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '37.24'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.41%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.065.27'
$ws.Range("E50").Value = '  -1.33%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1220'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +5.65%  '
